$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(", titular de la Cédula de Identidad N.º `${feligres_cedula}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${feligres_cedula_texto}", 2)
Write-Host "Found: $found"
